# Generate Report for Handback
#
# This sets the handback status on the Overview sheet, stamps the
# zh-cn / de-de detail sheets with "Latest Target File" / "Latest
# Handback File" links for the two rows that now have a handback, and
# records the handback timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Overview sheet: status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both localization columns
#    (zh-cn / de-de) of both rows.
# ---------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B2").Value = "Handed back: in sync with en-US"
$ovw.Range("C2").Value = "Handed back: in sync with en-US"
$ovw.Range("B3").Value = "Handed back: in sync with en-US"
$ovw.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# Helper data describing the two content files.
# ---------------------------------------------------------------
$mdUrlBase = "https://github.com/OpenLocalizationTest/oltest/blob/a26c2ce4ffc3aeae213fd111388237465e00a0c4/e2e"
$file1 = "bfc56466-f424-4c71-a9e5-4b645e843490"
$file2 = "eb5361da-e598-49e1-a781-2298c27002b2"

# ---------------------------------------------------------------
# 2. zh-cn sheet: add "Latest Target File" (F) / "Latest Handback
#    File" (G) hyperlinked cells on rows 2 and 3, then rebuild the
#    hyperlink collection in left-to-right / top-to-bottom order so
#    relationship ids line up the way a freshly generated report
#    would number them.
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("F2").Value = "$file1.md"
$wsZh.Range("G2").Value = "$file1.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf"
$wsZh.Range("F3").Value = "$file2.md"
$wsZh.Range("G3").Value = "$file2.7b12573b727babc04b54bfddf3d680bab711f64c.zh-cn.xlf"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$mdUrlBase/$file1.md", [Type]::Missing, [Type]::Missing, "$file1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e4c1973a36e434833cbfef11cc03e03f64d4723/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$file1.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "$file1.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "$mdUrlBase/$file1.md", [Type]::Missing, [Type]::Missing, "$file1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e4c1973a36e434833cbfef11cc03e03f64d4723/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$file1.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "$file1.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$mdUrlBase/$file2.md", [Type]::Missing, [Type]::Missing, "$file2.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e4c1973a36e434833cbfef11cc03e03f64d4723/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$file2.7b12573b727babc04b54bfddf3d680bab711f64c.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "$file2.7b12573b727babc04b54bfddf3d680bab711f64c.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "$mdUrlBase/$file2.md", [Type]::Missing, [Type]::Missing, "$file2.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e4c1973a36e434833cbfef11cc03e03f64d4723/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$file2.7b12573b727babc04b54bfddf3d680bab711f64c.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "$file2.7b12573b727babc04b54bfddf3d680bab711f64c.zh-cn.xlf")

# Restyle the newly-added link cells to match the workbook's existing
# hyperlink look (underlined, cornflower-blue) instead of Excel's
# built-in theme hyperlink style.
foreach ($addr in @("F2", "G2", "F3", "G3")) {
    $f = $wsZh.Range($addr).Font
    $f.Underline = $true
    $f.Color = 15570276
}

# ---------------------------------------------------------------
# 3. de-de sheet: same F/G additions, plus the handback timestamp
#    (column H) moves from the "never handed back" placeholder to the
#    actual handback datetime.
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("F2").Value = "$file1.md"
$wsDe.Range("G2").Value = "$file1.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf"
$wsDe.Range("H2").Value = "2016-03-24 00:50:47"
$wsDe.Range("F3").Value = "$file2.md"
$wsDe.Range("G3").Value = "$file2.7b12573b727babc04b54bfddf3d680bab711f64c.de-de.xlf"
$wsDe.Range("H3").Value = "2016-03-24 00:50:47"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$mdUrlBase/$file1.md", [Type]::Missing, [Type]::Missing, "$file1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6d843090e2ed1158392d1a82d7f435837cdb140e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$file1.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf", [Type]::Missing, [Type]::Missing, "$file1.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "$mdUrlBase/$file1.md", [Type]::Missing, [Type]::Missing, "$file1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6d843090e2ed1158392d1a82d7f435837cdb140e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$file1.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf", [Type]::Missing, [Type]::Missing, "$file1.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$mdUrlBase/$file2.md", [Type]::Missing, [Type]::Missing, "$file2.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6d843090e2ed1158392d1a82d7f435837cdb140e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$file2.7b12573b727babc04b54bfddf3d680bab711f64c.de-de.xlf", [Type]::Missing, [Type]::Missing, "$file2.7b12573b727babc04b54bfddf3d680bab711f64c.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "$mdUrlBase/$file2.md", [Type]::Missing, [Type]::Missing, "$file2.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6d843090e2ed1158392d1a82d7f435837cdb140e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$file2.7b12573b727babc04b54bfddf3d680bab711f64c.de-de.xlf", [Type]::Missing, [Type]::Missing, "$file2.7b12573b727babc04b54bfddf3d680bab711f64c.de-de.xlf")

foreach ($addr in @("F2", "G2", "F3", "G3")) {
    $f = $wsDe.Range($addr).Font
    $f.Underline = $true
    $f.Color = 15570276
}
